$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20 (shifts existing rows 20..138 down to 21..139)
$ws.Rows(20).Insert()

# Populate the new row 20 with the latest weekly price record, copying the
# constant columns from the row below (now row 21) and setting the new data.
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44670
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 100112044
$ws.Cells.Item(20, 7).Value = "Perejil"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 3000
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = 2250
$ws.Cells.Item(20, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(20, 16).Value = 1500
$ws.Cells.Item(20, 17).Value = 1.5
$ws.Cells.Item(20, 18).Value = "Hortaliza"
